# Update the NATMI TPM-derived values in the sheet to reflect the new TPM-based
# recalculation (commit: "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("M2").Value = 0.00535
$ws.Range("N2").Value = 0.01605
$ws.Range("O2").Value = 0.003591913026022235
$ws.Range("P2").Value = 0.003591913026022235
$ws.Range("Q2").Value = 0.0333845778
$ws.Range("R2").Value = 0.3004612001999999
$ws.Range("S2").Value = 0.00006222021381574401
$ws.Range("T2").Value = 0.00006222021381574401

$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("O3").Value = 0.9964080869739778
$ws.Range("P3").Value = 0.9964080869739778
$ws.Range("Q3").Value = 9.260987963555998
$ws.Range("R3").Value = 83.34889167200399
$ws.Range("S3").Value = 0.01726008502157802
$ws.Range("T3").Value = 0.01726008502157802

$ws.Range("I4").Value = 0.9592798330716089
$ws.Range("J4").Value = 0.9592798330716091
$ws.Range("M4").Value = 0.00535
$ws.Range("N4").Value = 0.01605
$ws.Range("O4").Value = 0.003591913026022235
$ws.Range("P4").Value = 0.003591913026022235
$ws.Range("Q4").Value = 1.84878119765
$ws.Range("R4").Value = 16.63903077885
$ws.Range("S4").Value = 0.003445649728010347
$ws.Range("T4").Value = 0.003445649728010348

$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("O5").Value = 0.9964080869739778
$ws.Range("P5").Value = 0.9964080869739778
$ws.Range("R5").Value = 4615.720009739576
$ws.Range("S5").Value = 0.9558341833435986
$ws.Range("T5").Value = 0.9558341833435988

$ws.Range("G6").Value = 8.428738666666666
$ws.Range("I6").Value = 0.02339786169299727
$ws.Range("J6").Value = 0.02339786169299728
$ws.Range("M6").Value = 0.00535
$ws.Range("N6").Value = 0.01605
$ws.Range("O6").Value = 0.003591913026022235
$ws.Range("P6").Value = 0.003591913026022235
$ws.Range("Q6").Value = 0.04509375186666666
$ws.Range("R6").Value = 0.4058437667999999
$ws.Range("S6").Value = 0.00008404308419614356
$ws.Range("T6").Value = 0.00008404308419614358

$ws.Range("G7").Value = 8.428738666666666
$ws.Range("I7").Value = 0.02339786169299727
$ws.Range("J7").Value = 0.02339786169299728
$ws.Range("O7").Value = 0.9964080869739778
$ws.Range("P7").Value = 0.9964080869739778
$ws.Range("Q7").Value = 12.50915005637066
$ws.Range("S7").Value = 0.02331381860880113
$ws.Range("T7").Value = 0.02331381860880113
